# Change highlight colors on Slide "Step 4" (Level1Scene-related Sprite.c slide)
# and Slide "Step 5" (Level1Scene.c slide) content placeholders.

$p = $ppt.ActivePresentation

# --- Slide 11: "Verify that both colored and textured meshes are still
#     drawn correctly" -> highlight changes from yellow (FFFF00) to green (00FF00)
$slide11 = $p.Slides.Item(11)
$shape11 = $slide11.Shapes.Item(2)
$para11 = $shape11.TextFrame.TextRange.Paragraphs(3)
$para11.Font.Highlight.RGB = 65280  # 00FF00 (green) in BGR-packed OLE color order

# --- Slide 12: "Change the Planet Entity from PlanetJump.txt to
#     PlanetBounce.txt" -> add a green (00FF00) highlight
$slide12 = $p.Slides.Item(12)
$shape12 = $slide12.Shapes.Item(2)
$para12 = $shape12.TextFrame.TextRange.Paragraphs(2)
$para12.Font.Highlight.RGB = 65280  # 00FF00 (green) in BGR-packed OLE color order
